# Weekly update: a new price-report entry for "Terminal La Palmera de La
# Serena - Berenjena" is inserted as the new row 247, pushing all the
# existing entries (old rows 247-287) down by one row (new rows 248-288)
# and adding one brand-new row at the bottom (288).
#
# This mirrors what happened in the source workbook: a fresh row was
# inserted above the old row 247, the older rows shifted down intact, and
# the brand-new row 247 was then filled in with its own data (date serial
# 45154, volume 380; the remaining fields happen to coincide with the
# values the old row 247 used to hold).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 247; existing rows 247:287 shift to 248:288.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new market entry.
$ws.Cells.Item(247, 1).Value = 8
$ws.Cells.Item(247, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(247, 3).Value = "Coquimbo"
$ws.Cells.Item(247, 4).Value = 45154
$ws.Cells.Item(247, 5).Value = 4
$ws.Cells.Item(247, 6).Value = 100112001
$ws.Cells.Item(247, 7).Value = "Berenjena"
$ws.Cells.Item(247, 8).Value = "Sin especificar"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 380
$ws.Cells.Item(247, 11).Value = 11000
$ws.Cells.Item(247, 12).Value = 12000
$ws.Cells.Item(247, 13).Value = 11500
$ws.Cells.Item(247, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(247, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(247, 16).Value = 230
$ws.Cells.Item(247, 17).Value = 50
$ws.Cells.Item(247, 18).Value = "Hortaliza"
